$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.005243333333333
$ws.Cells.Item(2, 8).Value = 3.01573
$ws.Cells.Item(2, 9).Value = 0.07224874268505826
$ws.Cells.Item(2, 10).Value = 0.07224874268505825
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 31.82741333333333
$ws.Cells.Item(2, 14).Value = 95.48223999999999
$ws.Cells.Item(2, 15).Value = 0.114390792932228
$ws.Cells.Item(2, 16).Value = 0.114390792932228
$ws.Cells.Item(2, 17).Value = 31.99429507057777
$ws.Cells.Item(2, 18).Value = 287.9486556351999
$ws.Cells.Item(2, 19).Value = 0.008264590964100323
$ws.Cells.Item(2, 20).Value = 0.008264590964100323

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.005243333333333
$ws.Cells.Item(3, 8).Value = 3.01573
$ws.Cells.Item(3, 9).Value = 0.07224874268505826
$ws.Cells.Item(3, 10).Value = 0.07224874268505825
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 85.46317833333335
$ws.Cells.Item(3, 14).Value = 256.389535
$ws.Cells.Item(3, 15).Value = 0.307162904935779
$ws.Cells.Item(3, 16).Value = 0.307162904935779
$ws.Cells.Item(3, 17).Value = 85.91129026506111
$ws.Cells.Item(3, 18).Value = 773.20161238555
$ws.Cells.Item(3, 19).Value = 0.02219213368110011
$ws.Cells.Item(3, 20).Value = 0.02219213368110011

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.005243333333333
$ws.Cells.Item(4, 8).Value = 3.01573
$ws.Cells.Item(4, 9).Value = 0.07224874268505826
$ws.Cells.Item(4, 10).Value = 0.07224874268505825
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 122.2478306666667
$ws.Cells.Item(4, 14).Value = 366.743492
$ws.Cells.Item(4, 15).Value = 0.4393704929064738
$ws.Cells.Item(4, 16).Value = 0.4393704929064738
$ws.Cells.Item(4, 17).Value = 122.8888167921289
$ws.Cells.Item(4, 18).Value = 1105.99935112916
$ws.Cells.Item(4, 19).Value = 0.03174396568540704
$ws.Cells.Item(4, 20).Value = 0.03174396568540704

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.005243333333333
$ws.Cells.Item(5, 8).Value = 3.01573
$ws.Cells.Item(5, 9).Value = 0.07224874268505826
$ws.Cells.Item(5, 10).Value = 0.07224874268505825
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 38.69562533333333
$ws.Cells.Item(5, 14).Value = 116.086876
$ws.Cells.Item(5, 15).Value = 0.1390758092255191
$ws.Cells.Item(5, 16).Value = 0.1390758092255191
$ws.Cells.Item(5, 17).Value = 38.89851939549777
$ws.Cells.Item(5, 18).Value = 350.0866745594799
$ws.Cells.Item(5, 19).Value = 0.01004805235445078
$ws.Cells.Item(5, 20).Value = 0.01004805235445078

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 10.25983933333333
$ws.Cells.Item(6, 8).Value = 30.779518
$ws.Cells.Item(6, 9).Value = 0.7373940889775011
$ws.Cells.Item(6, 10).Value = 0.737394088977501
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 31.82741333333333
$ws.Cells.Item(6, 14).Value = 95.48223999999999
$ws.Cells.Item(6, 15).Value = 0.114390792932228
$ws.Cells.Item(6, 16).Value = 0.114390792932228
$ws.Cells.Item(6, 17).Value = 326.5441471955911
$ws.Cells.Item(6, 18).Value = 2938.89732476032
$ws.Cells.Item(6, 19).Value = 0.08435109454167425
$ws.Cells.Item(6, 20).Value = 0.08435109454167426

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 10.25983933333333
$ws.Cells.Item(7, 8).Value = 30.779518
$ws.Cells.Item(7, 9).Value = 0.7373940889775011
$ws.Cells.Item(7, 10).Value = 0.737394088977501
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 85.46317833333335
$ws.Cells.Item(7, 14).Value = 256.389535
$ws.Cells.Item(7, 15).Value = 0.307162904935779
$ws.Cells.Item(7, 16).Value = 0.307162904935779
$ws.Cells.Item(7, 17).Value = 876.8384786160146
$ws.Cells.Item(7, 18).Value = 7891.54630754413
$ws.Cells.Item(7, 19).Value = 0.2265001104528016
$ws.Cells.Item(7, 20).Value = 0.2265001104528015

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 10.25983933333333
$ws.Cells.Item(8, 8).Value = 30.779518
$ws.Cells.Item(8, 9).Value = 0.7373940889775011
$ws.Cells.Item(8, 10).Value = 0.737394088977501
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 122.2478306666667
$ws.Cells.Item(8, 14).Value = 366.743492
$ws.Cells.Item(8, 15).Value = 0.4393704929064738
$ws.Cells.Item(8, 16).Value = 0.4393704929064738
$ws.Cells.Item(8, 17).Value = 1254.24310148854
$ws.Cells.Item(8, 18).Value = 11288.18791339686
$ws.Cells.Item(8, 19).Value = 0.3239892043403649
$ws.Cells.Item(8, 20).Value = 0.3239892043403649

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 10.25983933333333
$ws.Cells.Item(9, 8).Value = 30.779518
$ws.Cells.Item(9, 9).Value = 0.7373940889775011
$ws.Cells.Item(9, 10).Value = 0.737394088977501
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 38.69562533333333
$ws.Cells.Item(9, 14).Value = 116.086876
$ws.Cells.Item(9, 15).Value = 0.1390758092255191
$ws.Cells.Item(9, 16).Value = 0.1390758092255191
$ws.Cells.Item(9, 17).Value = 397.0108988228631
$ws.Cells.Item(9, 18).Value = 3573.098089405767
$ws.Cells.Item(9, 19).Value = 0.1025536796426604
$ws.Cells.Item(9, 20).Value = 0.1025536796426604

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.383875
$ws.Cells.Item(10, 8).Value = 1.151625
$ws.Cells.Item(10, 9).Value = 0.02758982345723265
$ws.Cells.Item(10, 10).Value = 0.02758982345723265
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 31.82741333333333
$ws.Cells.Item(10, 14).Value = 95.48223999999999
$ws.Cells.Item(10, 15).Value = 0.114390792932228
$ws.Cells.Item(10, 16).Value = 0.114390792932228
$ws.Cells.Item(10, 17).Value = 12.21774829333333
$ws.Cells.Item(10, 18).Value = 109.95973464
$ws.Cells.Item(10, 19).Value = 0.003156021782133028
$ws.Cells.Item(10, 20).Value = 0.003156021782133028

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.383875
$ws.Cells.Item(11, 8).Value = 1.151625
$ws.Cells.Item(11, 9).Value = 0.02758982345723265
$ws.Cells.Item(11, 10).Value = 0.02758982345723265
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 85.46317833333335
$ws.Cells.Item(11, 14).Value = 256.389535
$ws.Cells.Item(11, 15).Value = 0.307162904935779
$ws.Cells.Item(11, 16).Value = 0.307162904935779
$ws.Cells.Item(11, 17).Value = 32.80717758270834
$ws.Cells.Item(11, 18).Value = 295.2645982443751
$ws.Cells.Item(11, 19).Value = 0.00847457031978888
$ws.Cells.Item(11, 20).Value = 0.00847457031978888

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.383875
$ws.Cells.Item(12, 8).Value = 1.151625
$ws.Cells.Item(12, 9).Value = 0.02758982345723265
$ws.Cells.Item(12, 10).Value = 0.02758982345723265
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 122.2478306666667
$ws.Cells.Item(12, 14).Value = 366.743492
$ws.Cells.Item(12, 15).Value = 0.4393704929064738
$ws.Cells.Item(12, 16).Value = 0.4393704929064738
$ws.Cells.Item(12, 17).Value = 46.92788599716667
$ws.Cells.Item(12, 18).Value = 422.3509739745001
$ws.Cells.Item(12, 19).Value = 0.0121221543316069
$ws.Cells.Item(12, 20).Value = 0.0121221543316069

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.383875
$ws.Cells.Item(13, 8).Value = 1.151625
$ws.Cells.Item(13, 9).Value = 0.02758982345723265
$ws.Cells.Item(13, 10).Value = 0.02758982345723265
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 38.69562533333333
$ws.Cells.Item(13, 14).Value = 116.086876
$ws.Cells.Item(13, 15).Value = 0.1390758092255191
$ws.Cells.Item(13, 16).Value = 0.1390758092255191
$ws.Cells.Item(13, 17).Value = 14.85428317483333
$ws.Cells.Item(13, 18).Value = 133.6885485735
$ws.Cells.Item(13, 19).Value = 0.003837077023703841
$ws.Cells.Item(13, 20).Value = 0.003837077023703841

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 2.264687
$ws.Cells.Item(14, 8).Value = 6.794061
$ws.Cells.Item(14, 9).Value = 0.162767344880208
$ws.Cells.Item(14, 10).Value = 0.162767344880208
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 31.82741333333333
$ws.Cells.Item(14, 14).Value = 95.48223999999999
$ws.Cells.Item(14, 15).Value = 0.114390792932228
$ws.Cells.Item(14, 16).Value = 0.114390792932228
$ws.Cells.Item(14, 17).Value = 72.07912921962665
$ws.Cells.Item(14, 18).Value = 648.71216297664
$ws.Cells.Item(14, 19).Value = 0.01861908564432042
$ws.Cells.Item(14, 20).Value = 0.01861908564432042

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 2.264687
$ws.Cells.Item(15, 8).Value = 6.794061
$ws.Cells.Item(15, 9).Value = 0.162767344880208
$ws.Cells.Item(15, 10).Value = 0.162767344880208
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 85.46317833333335
$ws.Cells.Item(15, 14).Value = 256.389535
$ws.Cells.Item(15, 15).Value = 0.307162904935779
$ws.Cells.Item(15, 16).Value = 0.307162904935779
$ws.Cells.Item(15, 17).Value = 193.5473489501817
$ws.Cells.Item(15, 18).Value = 1741.926140551635
$ws.Cells.Item(15, 19).Value = 0.04999609048208849
$ws.Cells.Item(15, 20).Value = 0.04999609048208849

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 2.264687
$ws.Cells.Item(16, 8).Value = 6.794061
$ws.Cells.Item(16, 9).Value = 0.162767344880208
$ws.Cells.Item(16, 10).Value = 0.162767344880208
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 122.2478306666667
$ws.Cells.Item(16, 14).Value = 366.743492
$ws.Cells.Item(16, 15).Value = 0.4393704929064738
$ws.Cells.Item(16, 16).Value = 0.4393704929064738
$ws.Cells.Item(16, 17).Value = 276.8530728890013
$ws.Cells.Item(16, 18).Value = 2491.677656001012
$ws.Cells.Item(16, 19).Value = 0.071515168549095
$ws.Cells.Item(16, 20).Value = 0.07151516854909501

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 2.264687
$ws.Cells.Item(17, 8).Value = 6.794061
$ws.Cells.Item(17, 9).Value = 0.162767344880208
$ws.Cells.Item(17, 10).Value = 0.162767344880208
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 38.69562533333333
$ws.Cells.Item(17, 14).Value = 116.086876
$ws.Cells.Item(17, 15).Value = 0.1390758092255191
$ws.Cells.Item(17, 16).Value = 0.1390758092255191
$ws.Cells.Item(17, 17).Value = 87.63347964927065
$ws.Cells.Item(17, 18).Value = 788.701316843436
$ws.Cells.Item(17, 19).Value = 0.02263700020470408
$ws.Cells.Item(17, 20).Value = 0.02263700020470408
